# Update data: 7 May 2021
# Adds the newest (April 2021) unemployment data point to both the
# "Canada" sheet (sheet1) and the "Province" sheet (sheet2).

$wb = $excel.ActiveWorkbook

$wsCanada = $wb.Worksheets.Item("Canada")
$wsProvince = $wb.Worksheets.Item("Province")

# ---------------------------------------------------------------
# Sheet "Canada": append row 17
# ---------------------------------------------------------------
$newDate = 44287

$wsCanada.Range("A17").Value2 = $newDate
$wsCanada.Range("A17").NumberFormat = $wsCanada.Range("A16").NumberFormat

$wsCanada.Range("B17").Value2 = "Canada"
$wsCanada.Range("B17").NumberFormat = $wsCanada.Range("B16").NumberFormat

$wsCanada.Range("D17").Value2 = 1640.3
$wsCanada.Range("E17").Value2 = 1166.9000000000001

$wsCanada.Range("C17").Formula = "=(D17-E17)/E17*100"

# ---------------------------------------------------------------
# Sheet "Province": append rows 152-161 (one per province, in the
# same order used throughout the sheet)
# ---------------------------------------------------------------
$provinceRows = @(
    @{ Row = 152; Name = "Newfoundland & Labrador"; D = 35.299999999999997;  E = 32.299999999999997 },
    @{ Row = 153; Name = "Prince Edward Island";     D = 7;                   E = 7.5 },
    @{ Row = 154; Name = "Nova Scotia";               D = 41.2;               E = 34.6 },
    @{ Row = 155; Name = "New Brunswick";             D = 33.9;               E = 32.299999999999997 },
    @{ Row = 156; Name = "Quebec";                    D = 296.39999999999998; E = 226 },
    @{ Row = 157; Name = "Ontario";                   D = 716.8;              E = 468.8 },
    @{ Row = 158; Name = "Manitoba";                  D = 52.7;               E = 36.299999999999997 },
    @{ Row = 159; Name = "Saskatchewan";               D = 39.4;               E = 34.200000000000003 },
    @{ Row = 160; Name = "Alberta";                   D = 218.8;              E = 164.4 },
    @{ Row = 161; Name = "British Columbia";          D = 198.9;              E = 130.5 }
)

foreach ($item in $provinceRows) {
    $r = $item.Row

    $wsProvince.Range("A$r").Value2 = $newDate
    $wsProvince.Range("A$r").NumberFormat = $wsProvince.Range("A151").NumberFormat

    $wsProvince.Range("B$r").Value2 = $item.Name
    if ($item.Row -eq 152) {
        $wsProvince.Range("B$r").NumberFormat = $wsProvince.Range("A151").NumberFormat
    }

    $wsProvince.Range("D$r").Value2 = $item.D
    $wsProvince.Range("E$r").Value2 = $item.E

    $wsProvince.Range("C$r").Formula = "=(D" + $r + "-E" + $r + ")/E" + $r + "*100"
}

# ---------------------------------------------------------------
# Restore view/selection state to match the edited workbook:
#  - Canada sheet: selection C16:C17, scrolled to show row 4
#  - Province sheet: selection D162, scrolled to show row 151,
#    and keep it the active/tab-selected sheet
# ---------------------------------------------------------------
$wsCanada.Activate()
$winCanada = $excel.ActiveWindow
$winCanada.ScrollColumn = 1
$winCanada.ScrollRow = 4
$wsCanada.Range("C16:C17").Select()

$wsProvince.Activate()
$winProvince = $excel.ActiveWindow
$winProvince.ScrollColumn = 1
$winProvince.ScrollRow = 151
$wsProvince.Range("D162").Select()
